$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 180.14285
$ws.Range("I9").Value = 208.81818
$ws.Range("J9").Value = 75
$ws.Range("K9").Value = 208.81818
$ws.Range("L9").Value = 75
$ws.Range("M9").Value = -39.81818000000001
$ws.Range("N9").Value = -413

$ws.Range("H33").Value = 263.16666
$ws.Range("I33").Value = 215.8
$ws.Range("K33").Value = 215.8
$ws.Range("M33").Value = 13.19999999999999

$ws.Range("H97").Value = 2655.5715
$ws.Range("J97").Value = 2655.5715
$ws.Range("L97").Value = 7966.7145
$ws.Range("N97").Value = -8958.7145

$ws.Range("H111").Value = 899
$ws.Range("I111").Value = 931
$ws.Range("J111").Value = 739
$ws.Range("K111").Value = 2793
$ws.Range("L111").Value = 2217
$ws.Range("M111").Value = 274
$ws.Range("N111").Value = -8351

$ws.Range("H112").Value = 2518.913
$ws.Range("J112").Value = 2823.2
$ws.Range("L112").Value = 8469.599999999999
$ws.Range("N112").Value = -10685.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 660.25
$ws.Range("I2").Value = 660.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 660.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -547.25
$ws.Range("N2").ClearContents()

$ws.Range("H5").Value = 3831.6667
$ws.Range("I5").Value = 998.25
$ws.Range("K5").Value = 998.25
$ws.Range("M5").Value = -886.25

$ws.Range("H32").Value = 7274.7915
$ws.Range("I32").Value = 6480.7144
$ws.Range("K32").Value = 6480.7144
$ws.Range("M32").Value = -6193.7144

$ws.Range("H61").Value = 2478.6538
$ws.Range("I61").Value = 2214.0952
$ws.Range("J61").Value = 3589.8
$ws.Range("K61").Value = 2214.0952
$ws.Range("L61").Value = 3589.8
$ws.Range("M61").Value = -2002.0952
$ws.Range("N61").Value = -4013.8

$ws.Range("H116").Value = 660.25
$ws.Range("I116").Value = 660.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 660.25
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1633.75
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 2086171.6
$ws.Range("I132").Value = 3127595.2
$ws.Range("J132").Value = 3324.25
$ws.Range("K132").Value = 9382785.600000001
$ws.Range("L132").Value = 9972.75
$ws.Range("M132").Value = -9380255.600000001
$ws.Range("N132").Value = -15032.75

$ws.Range("H136").Value = 2478.6538
$ws.Range("I136").Value = 2214.0952
$ws.Range("J136").Value = 3589.8
$ws.Range("K136").Value = 6642.285600000001
$ws.Range("L136").Value = 10769.4
$ws.Range("M136").Value = -4092.285600000001
$ws.Range("N136").Value = -15869.4

$ws.Range("H141").Value = 53143
$ws.Range("J141").Value = 53143
$ws.Range("L141").Value = 53143
$ws.Range("N141").Value = -63503

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 660.25
$ws.Range("I3").Value = 660.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 660.25
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -546.25
$ws.Range("N3").ClearContents()

$ws.Range("H4").Value = 3831.6667
$ws.Range("I4").Value = 998.25
$ws.Range("K4").Value = 998.25
$ws.Range("M4").Value = -883.25

$ws.Range("H140").Value = 84592.664
$ws.Range("J140").Value = 84592.664
$ws.Range("L140").Value = 84592.664
$ws.Range("N140").Value = -94952.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 136
$ws.Range("I7").Value = 104.8125
$ws.Range("K7").Value = 104.8125
$ws.Range("M7").Value = 8.1875

$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -330
$ws.Range("N14").ClearContents()

$ws.Range("H107").Value = 554.5217
$ws.Range("J107").Value = 705
$ws.Range("L107").Value = 705
$ws.Range("N107").Value = -4545

$ws.Range("H134").Value = 3476.75
$ws.Range("J134").Value = 5200
$ws.Range("L134").Value = 15600
$ws.Range("N134").Value = -20670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 728.38464
$ws.Range("I6").Value = 364.8889
$ws.Range("J6").Value = 1546.25
$ws.Range("K6").Value = 1094.6667
$ws.Range("L6").Value = 4638.75
$ws.Range("M6").Value = -981.6667
$ws.Range("N6").Value = -4864.75

$ws.Range("H11").Value = 58823748
$ws.Range("I11").Value = 160.3
$ws.Range("K11").Value = 480.9
$ws.Range("M11").Value = -340.9

$ws.Range("H38").Value = 160
$ws.Range("I38").Value = 161.66667
$ws.Range("J38").Value = 155
$ws.Range("K38").Value = 485.00001
$ws.Range("L38").Value = 465
$ws.Range("M38").Value = -138.00001
$ws.Range("N38").Value = -1159

$ws.Range("H137").Value = 7663.8
$ws.Range("I137").Value = 2579.7144
$ws.Range("K137").Value = 7739.1432
$ws.Range("M137").Value = -2639.1432

$ws.Range("H140").Value = 1273.32
$ws.Range("I140").Value = 1118.0416
$ws.Range("K140").Value = 3354.1248
$ws.Range("M140").Value = 1825.8752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 8001.3335
$ws.Range("I18").Value = 2002
$ws.Range("K18").Value = 2002
$ws.Range("M18").Value = -1709

$ws.Range("H132").Value = 12583.088
$ws.Range("I132").Value = 15243.615
$ws.Range("K132").Value = 45730.845
$ws.Range("M132").Value = -43200.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3492.7334
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 3914.6924
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 3914.6924
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -4504.6924

$ws.Range("H27").Value = 3492.7334
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 3914.6924
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 3914.6924
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -4128.6924

$ws.Range("H74").Value = 49661.332
$ws.Range("I74").Value = 41992
$ws.Range("K74").Value = 41992
$ws.Range("M74").Value = -40994

$ws.Range("H77").Value = 49661.332
$ws.Range("I77").Value = 41992
$ws.Range("K77").Value = 125976
$ws.Range("M77").Value = -120984

$ws.Range("H82").Value = 1495
$ws.Range("I82").Value = 1742.5
$ws.Range("K82").Value = 1742.5
$ws.Range("M82").Value = -1381.5

$ws.Range("H85").Value = 1495
$ws.Range("I85").Value = 1742.5
$ws.Range("K85").Value = 1742.5
$ws.Range("M85").Value = -494.5

$ws.Range("H100").Value = 4061.6
$ws.Range("I100").Value = 3827
$ws.Range("K100").Value = 3827
$ws.Range("M100").Value = -3286

$ws.Range("H122").Value = 4158
$ws.Range("I122").Value = 3767.25
$ws.Range("K122").Value = 11301.75
$ws.Range("M122").Value = -8851.75

$ws.Range("H132").Value = 2858.3076
$ws.Range("J132").Value = 2878
$ws.Range("L132").Value = 8634
$ws.Range("N132").Value = -13694

$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws.Range("H136").Value = 7206.16
$ws.Range("I136").Value = 3938.3333
$ws.Range("K136").Value = 11814.9999
$ws.Range("M136").Value = -9264.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1850.2222
$ws.Range("I96").Value = 1456.5
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 1456.5
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -83.5
$ws.Range("N96").Value = -7746

$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H100").Value = 1132.36
$ws.Range("I100").Value = 1171.8636
$ws.Range("J100").Value = 842.6667
$ws.Range("K100").Value = 2343.7272
$ws.Range("L100").Value = 1685.3334
$ws.Range("M100").Value = -1802.7272
$ws.Range("N100").Value = -2767.3334

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 50003
$ws.Range("J105").Value = 50003
$ws.Range("L105").Value = 50003
$ws.Range("N105").Value = -56991

$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H132").Value = 3215.8
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3215.8
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9647.400000000001
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14707.4

$ws.Range("H136").Value = 203676.23
$ws.Range("I136").Value = 213745.16
$ws.Range("K136").Value = 641235.48
$ws.Range("M136").Value = -638685.48
